$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Latest recorded run for each test user: token (column C) + id (column D).
# Row 2 = daniel5f, Row 3 = Jorge2525, Row 4 = mario35.

$ws.Range("C2").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6ImRhbmllbDVmIiwicGFzc3dvcmQiOiJBejI1Mjg4QCIsImlhdCI6MTcwNDE0NjA4OH0.MsNLx9d4FjW5892u7y2I7dd6kkk1SERNZX3cxAzpyfg"
$ws.Range("D2").Value = "6319d234-7d77-4577-9d84-f66e4d1389db"

$ws.Range("C3").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6IkpvcmdlMjUyNSIsInBhc3N3b3JkIjoiYXNUMzU2NDQ0QCIsImlhdCI6MTcwNDE0NjA4OX0.XnIo44G9clbwFbsl4bbiEUpZqOqOeWF9EQ4a51Hw0O0"
$ws.Range("D3").Value = "1696c4b7-8035-488e-9e85-fdd83bd48c08"

$ws.Range("C4").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6Im1hcmlvMzUiLCJwYXNzd29yZCI6Im1BcmlvdXVnQDMiLCJpYXQiOjE3MDQxNDYwOTB9.dsDutbtbyTMgHbFLLCEzJ2xN0MCbRV3QuBzpZqSm8Xw"
$ws.Range("D4").Value = "0ea114b4-4e9f-4ecf-bb57-5b9773769adf"
